$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.76"
$ws.Range("D3").Value = "'23.04"
$ws.Range("D4").Value = "'6.396"
$ws.Range("D5").Value = "'0.06271"
$ws.Range("D7").Value = "'6.667"
$ws.Range("D8").Value = "'1.360"
$ws.Range("D9").Value = "'0.8308"
$ws.Range("D10").Value = "'0.01373"
$ws.Range("D11").Value = "'0.1628"
$ws.Range("D12").Value = "'0.08313"
$ws.Range("D13").Value = "'0.03436"
$ws.Range("D14").Value = "'0.03088"
$ws.Range("D15").Value = "'0.09315"
$ws.Range("D16").Value = "'3.874"
$ws.Range("D17").Value = "'0.001654"
$ws.Range("D18").Value = "'0.04776"
$ws.Range("D19").Value = "'0.006303"
$ws.Range("D20").Value = "'0.005681"
$ws.Range("D21").Value = "'0.001092"
$ws.Range("D23").Value = "'3.713"
$ws.Range("D27").Value = "'0.0002680"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("D41").Value = "'0.007056"
$ws.Range("D43").Value = "'0.003700"
$ws.Range("D44").Value = "'0.01220"
$ws.Range("D45").Value = "'0.00006269"
$ws.Range("D48").Value = "'0.7701"
$ws.Range("D49").Value = "'0.02058"
$ws.Range("D50").Value = "'0.00002300"
